{"js": "// Replace each two-digit-by-two-digit multiplication prompt in the\n// worksheet table with its new pair of operands, preserving all\n// paragraph/run formatting (font, size, alignment, etc.).\nconst replacements = [\n  [\"89\u00d770=\", \"18\u00d750=\"],\n  [\"60\u00d779=\", \"45\u00d795=\"],\n  [\"68\u00d721=\", \"78\u00d790=\"],\n  [\"77\u00d743=\", \"47\u00d754=\"],\n  [\"78\u00d796=\", \"34\u00d738=\"],\n  [\"15\u00d757=\", \"22\u00d747=\"],\n  [\"79\u00d734=\", \"98\u00d788=\"],\n  [\"26\u00d741=\", \"58\u00d779=\"],\n  [\"19\u00d718=\", \"47\u00d765=\"],\n  [\"82\u00d717=\", \"11\u00d732=\"],\n  [\"20\u00d773=\", \"96\u00d727=\"],\n  [\"70\u00d736=\", \"89\u00d758=\"],\n  [\"43\u00d741=\", \"98\u00d745=\"],\n  [\"19\u00d735=\", \"11\u00d774=\"],\n  [\"49\u00d778=\", \"24\u00d727=\"],\n  [\"21\u00d744=\", \"37\u00d774=\"],\n  [\"64\u00d724=\", \"67\u00d743=\"],\n  [\"86\u00d744=\", \"57\u00d740=\"],\n  [\"44\u00d751=\", \"96\u00d712=\"],\n  [\"88\u00d773=\", \"54\u00d779=\"],\n  [\"36\u00d797=\", \"26\u00d777=\"],\n  [\"14\u00d725=\", \"34\u00d763=\"],\n  [\"98\u00d782=\", \"55\u00d776=\"],\n  [\"90\u00d784=\", \"93\u00d783=\"],\n  [\"92\u00d794=\", \"29\u00d726=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  // Every prompt string is unique in the document, so matchCase search\n  // with no wildcards is sufficient to land on exactly the right cell.\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    // Replacing the matched range's text in place keeps the run's\n    // formatting (font, size, paragraph alignment, etc.) untouched.\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-by-two-digit multiplication prompt in the\n# worksheet table with its new pair of operands, preserving all\n# paragraph/run formatting (font, size, alignment, etc.). Each source\n# prompt is unique in the document, so a simple Find/Replace per pair\n# (scoped to the whole document body) lands on the correct cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"89\u00d770=\", \"18\u00d750=\"),\n    @(\"60\u00d779=\", \"45\u00d795=\"),\n    @(\"68\u00d721=\", \"78\u00d790=\"),\n    @(\"77\u00d743=\", \"47\u00d754=\"),\n    @(\"78\u00d796=\", \"34\u00d738=\"),\n    @(\"15\u00d757=\", \"22\u00d747=\"),\n    @(\"79\u00d734=\", \"98\u00d788=\"),\n    @(\"26\u00d741=\", \"58\u00d779=\"),\n    @(\"19\u00d718=\", \"47\u00d765=\"),\n    @(\"82\u00d717=\", \"11\u00d732=\"),\n    @(\"20\u00d773=\", \"96\u00d727=\"),\n    @(\"70\u00d736=\", \"89\u00d758=\"),\n    @(\"43\u00d741=\", \"98\u00d745=\"),\n    @(\"19\u00d735=\", \"11\u00d774=\"),\n    @(\"49\u00d778=\", \"24\u00d727=\"),\n    @(\"21\u00d744=\", \"37\u00d774=\"),\n    @(\"64\u00d724=\", \"67\u00d743=\"),\n    @(\"86\u00d744=\", \"57\u00d740=\"),\n    @(\"44\u00d751=\", \"96\u00d712=\"),\n    @(\"88\u00d773=\", \"54\u00d779=\"),\n    @(\"36\u00d797=\", \"26\u00d777=\"),\n    @(\"14\u00d725=\", \"34\u00d763=\"),\n    @(\"98\u00d782=\", \"55\u00d776=\"),\n    @(\"90\u00d784=\", \"93\u00d783=\"),\n    @(\"92\u00d794=\", \"29\u00d726=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
